$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.529.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.300.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.584'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.41%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.634'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0976'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.81%  '

$ws.Range("E12").Value = '  +1.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.812.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.292.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.89%  '

$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '58.176.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("E21").Value = '  +6.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '302.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.69%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.20%  '

$ws.Range("E30").Value = '  -1.76%  '

$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.113'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.98%  '

$ws.Range("E34").Value = '  +13.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0526'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.26%  '

$ws.Range("E36").Value = '  +0.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.95%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '137.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("E42").Value = '  +2.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.88'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.90%  '

$ws.Range("E44").Value = '  -1.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.281'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.169.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.73%  '
